$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.965.32"
$ws.Cells.Item(2, 5).Value = "  -0.60%  "

$ws.Cells.Item(3, 4).Value = "1.562.98"
$ws.Cells.Item(3, 5).Value = "  -0.50%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "207.49"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.38%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.490"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.39%  "

$ws.Cells.Item(7, 5).Value = "  +0.26%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.09"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.59%  "

$ws.Cells.Item(9, 5).Value = "  -0.54%  "

$ws.Cells.Item(10, 5).Value = "  +1.52%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0858"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.03%  "

$ws.Cells.Item(12, 4).Value = "1.781.79"
$ws.Cells.Item(12, 5).Value = "  -0.04%  "

$ws.Cells.Item(13, 4).Value = "1.558.80"
$ws.Cells.Item(13, 5).Value = "  -1.01%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.76"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.39%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.520"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.97%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "62.12"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.08%  "

$ws.Cells.Item(17, 4).Value = "26.951.09"
$ws.Cells.Item(17, 5).Value = "  -0.52%  "

$ws.Cells.Item(18, 4).Value = "0.0₃0706"
$ws.Cells.Item(18, 5).Value = "  +0.99%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "217.18"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.99%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.36"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.07%  "

$ws.Cells.Item(21, 5).Value = "  +0.01%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.11"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.40%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.23"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.93%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.17%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "152.71"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.02%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.61"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.66%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.09"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.48%  "

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.105"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.86%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.11%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0471"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.35%  "

$ws.Cells.Item(31, 5).Value = "  +1.01%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.24"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.25%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.12"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.49%  "

$ws.Cells.Item(34, 4).Value = "1.421.00"
$ws.Cells.Item(34, 5).Value = "  -2.85%  "

$ws.Cells.Item(35, 5).Value = "  +2.13%  "

$ws.Cells.Item(36, 5).Value = "  +9.91%  "

$ws.Cells.Item(37, 5).Value = "  +1.93%  "

$ws.Cells.Item(38, 5).Value = "  +0.04%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.533"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.33%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.809"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.23%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.75"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.15%  "

$ws.Cells.Item(42, 5).Value = "  +0.18%  "

$ws.Cells.Item(43, 5).Value = "  +1.90%  "

$ws.Cells.Item(44, 5).Value = "  +2.05%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "64.85"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.26%  "

$ws.Cells.Item(46, 5).Value = "  -1.67%  "

$ws.Cells.Item(47, 4).Value = "1.696.41"
$ws.Cells.Item(47, 5).Value = "  +0.15%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "87.55"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.12%  "

$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.0₆0102"
$ws.Cells.Item(49, 5).Value = "  +8.42%  "

$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0521"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.68%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0958"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.05%  "

